$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.423.77"
$ws.Range("E2").Value = "  +2.25%  "
$ws.Range("D3").Value = "2.599.54"
$ws.Range("E3").Value = "  +1.89%  "
$ws.Range("E4").Value = "  +0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "586.31"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +6.92%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "142.94"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.20%  "
$ws.Range("E7").Value = "  -0.20%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.599"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.75%  "
$ws.Range("D9").Value = "2.610.58"
$ws.Range("E9").Value = "  +2.25%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.52"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.20%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.106"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.94%  "
$ws.Range("E12").Value = "  -3.16%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.373"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +6.23%  "
$ws.Range("D14").Value = "3.068.91"
$ws.Range("E14").Value = "  +2.29%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "24.75"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +7.76%  "
$ws.Range("D16").Value = "60.442.06"
$ws.Range("E16").Value = "  +2.26%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.0000142"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +4.27%  "
$ws.Range("D18").Value = "2.615.52"
$ws.Range("E18").Value = "  +2.14%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.37"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +11.45%  "
$ws.Range("E20").Value = "  +3.52%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "348.11"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +3.68%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.93"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +8.47%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("E24").Value = "  +10.31%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "62.95"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("E27").Value = "  +0.95%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.96"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +8.65%  "
$ws.Range("E29").Value = "  +5.41%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.88"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +12.94%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +4.61%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "163.97"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +3.72%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "19.49"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.03%  "
$ws.Range("E35").Value = "  +5.05%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.987"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +11.56%  "
$ws.Range("E37").Value = "  +7.24%  "
$ws.Range("E38").Value = "  +11.36%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "37.92"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.83%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.90"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +7.35%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "311.77"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +10.33%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.843"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.65%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "135.94"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.96%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "5.09"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +14.49%  "
$ws.Range("E45").Value = "  +2.70%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.995"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.33%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "19.79"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +6.57%  "
$ws.Range("E48").Value = "  +4.88%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.603"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +2.76%  "
$ws.Range("E50").Value = "  +9.69%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0242"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +4.34%  "
